$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new tag/response row at the bottom of the table
$ws.Range("A43").Value = "ถูกหนึ่ง"
$ws.Range("B43").Value = "ข้อ 2 ตอนเที่ยงวันพระจันทร์ไปไหน"

# Update selection to match the author's final cursor position
$ws.Range("B43").Select()
